$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the shared string used by column E row 2 (" L/17" -> "L/17"); this
# string is reused by every duplicated block below.
$ws.Range("E2").Value = "L/17"

# Duplicate the 5 data rows (2-6) three more times below the existing
# data, preserving formatting/styles by using Range.Copy with a
# destination range (not PasteSpecial, which drops cell styles here).
$ws.Range("A2:K6").Copy($ws.Range("A7"))
$ws.Range("A2:K6").Copy($ws.Range("A12"))
$ws.Range("A2:K6").Copy($ws.Range("A17"))

# Renumber column A (the running id) continuously for the new rows.
for ($i = 0; $i -lt 15; $i++) {
    $row = 7 + $i
    $ws.Cells.Item($row, 1).Value = 6 + $i
}

# Update the used range and the remembered selection to match the grown
# sheet.
$ws.Range("A1:K21").Select()
$ws.Range("C27").Select()
